$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D24","D30","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "66.256.34"
$ws.Range("E2").Value = "  -1.35%  "

# Row 3
$ws.Range("D3").Value = "3.458.50"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "581.80"
$ws.Range("E5").Value = "  -2.02%  "

# Row 6
$ws.Range("D6").Value = "174.73"
$ws.Range("E6").Value = "  -2.64%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  -2.37%  "

# Row 9
$ws.Range("D9").Value = "3.448.42"
$ws.Range("E9").Value = "  -0.62%  "

# Row 10
$ws.Range("E10").Value = "  -4.76%  "

# Row 11
$ws.Range("D11").Value = "6.87"
$ws.Range("E11").Value = "  -1.26%  "

# Row 12
$ws.Range("D12").Value = "0.419"
$ws.Range("E12").Value = "  -3.09%  "

# Row 13
$ws.Range("D13").Value = "4.040.79"
$ws.Range("E13").Value = "  -0.73%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "30.14"
$ws.Range("E14").Value = "  -4.83%  "

# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.132"
$ws.Range("E15").Value = "  -1.37%  "

# Row 16
$ws.Range("D16").Value = "66.157.92"
$ws.Range("E16").Value = "  -1.47%  "

# Row 17
$ws.Range("D17").Value = "0.0000172"
$ws.Range("E17").Value = "  -2.68%  "

# Row 18
$ws.Range("D18").Value = "3.461.62"
$ws.Range("E18").Value = "  -0.27%  "

# Row 19
$ws.Range("D19").Value = "5.96"
$ws.Range("E19").Value = "  -4.59%  "

# Row 20
$ws.Range("D20").Value = "13.77"
$ws.Range("E20").Value = "  -2.66%  "

# Row 21
$ws.Range("D21").Value = "378.64"
$ws.Range("E21").Value = "  -2.26%  "

# Row 22
$ws.Range("D22").Value = "7.79"
$ws.Range("E22").Value = "  -1.54%  "

# Row 23
$ws.Range("D23").Value = "0.547"
$ws.Range("E23").Value = "  +1.89%  "

# Row 24
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").Value = "71.93"
$ws.Range("E25").Value = "  -0.93%  "

# Row 26
$ws.Range("D26").Value = "5.71"
$ws.Range("E26").Value = "  -0.78%  "

# Row 27
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  -2.74%  "

# Row 28
$ws.Range("D28").Value = "9.78"
$ws.Range("E28").Value = "  -4.51%  "

# Row 29
$ws.Range("E29").Value = "  -0.65%  "

# Row 30
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.11%  "

# Row 31
$ws.Range("D31").Value = "24.15"
$ws.Range("E31").Value = "  +3.08%  "

# Row 32
$ws.Range("D32").Value = "5.84"
$ws.Range("E32").Value = "  -5.41%  "

# Row 33
$ws.Range("D33").Value = "1.99"
$ws.Range("E33").Value = "  -3.26%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.32"
$ws.Range("E34").Value = "  -6.71%  "

# Row 35
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("D36").Value = "7.11"
$ws.Range("E36").Value = "  -2.88%  "

# Row 37
$ws.Range("D37").Value = "1.57"
$ws.Range("E37").Value = "  -1.80%  "

# Row 38
$ws.Range("D38").Value = "159.68"
$ws.Range("E38").Value = "  -1.22%  "

# Row 39
$ws.Range("D39").Value = "29.38"
$ws.Range("E39").Value = "  +12.61%  "

# Row 40
$ws.Range("D40").Value = "0.885"
$ws.Range("E40").Value = "  +0.47%  "

# Row 41
$ws.Range("D41").Value = "1.77"
$ws.Range("E41").Value = "  -5.11%  "

# Row 42
$ws.Range("E42").Value = "  -2.63%  "

# Row 43
$ws.Range("D43").Value = "2.54"
$ws.Range("E43").Value = "  -9.57%  "

# Row 44
$ws.Range("D44").Value = "6.34"
$ws.Range("E44").Value = "  -6.76%  "

# Row 45
$ws.Range("D45").Value = "0.0694"
$ws.Range("E45").Value = "  -3.72%  "

# Row 46
$ws.Range("D46").Value = "2.664.69"
$ws.Range("E46").Value = "  -4.82%  "

# Row 47
$ws.Range("D47").Value = "40.29"
$ws.Range("E47").Value = "  -1.99%  "

# Row 48
$ws.Range("D48").Value = "24.46"
$ws.Range("E48").Value = "  -8.24%  "

# Row 49
$ws.Range("D49").Value = "0.0291"
$ws.Range("E49").Value = "  -2.48%  "

# Row 50
$ws.Range("D50").Value = "312.54"
$ws.Range("E50").Value = "  -5.96%  "

# Row 51
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -4.18%  "
